$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns, plus a few coin name/link
# changes (rows 21-22 swap, rows 50-51 shift) to match the refreshed feed.

$ws.Range("D2").Value = "25.544.69"
$ws.Range("E2").Value = "  +1.45%  "
$ws.Range("D3").Value = "1.670.97"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").Value = "'0.9966"
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").Value = "'237.38"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'0.9978"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "'0.4808"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "'0.2630"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").Value = "'0.06164"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").Value = "'0.07092"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("D11").Value = "1.662.90"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "'14.87"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").Value = "'0.6027"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").Value = "'4.419"
$ws.Range("E14").Value = "  -3.79%  "
$ws.Range("D15").Value = "'74.62"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").Value = "'0.9982"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "'0.9968"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "25.520.28"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "'0.000006808"
$ws.Range("E19").Value = "  +3.78%  "
$ws.Range("D20").Value = "'11.48"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'4.477"
$ws.Range("E21").Value = "  +1.22%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "1.872.89"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "'8.694"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "'5.380"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").Value = "'134.34"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "'15.11"
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("D28").Value = "'104.66"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").Value = "'1.701"
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").Value = "'3.972"
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").Value = "'3.684"
$ws.Range("E31").Value = "  +4.08%  "
$ws.Range("D32").Value = "'0.07677"
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("D33").Value = "'0.04365"
$ws.Range("E33").Value = "  -5.33%  "
$ws.Range("D34").Value = "'0.9977"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("D35").Value = "'2.610"
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("D36").Value = "'0.6226"
$ws.Range("E36").Value = "  +6.71%  "
$ws.Range("D37").Value = "'0.9508"
$ws.Range("E37").Value = "  +0.81%  "
$ws.Range("D38").Value = "'2.620"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "'0.8615"
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("D40").Value = "'0.9970"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "'0.01505"
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").Value = "'1.866"
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("D43").Value = "'98.02"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").Value = "'0.3791"
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("D45").Value = "'4.673"
$ws.Range("E45").Value = "  -4.05%  "
$ws.Range("D46").Value = "'0.1121"
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "'6.235"
$ws.Range("E47").Value = "  +2.93%  "
$ws.Range("D48").Value = "'0.05254"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").Value = "'29.59"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.377"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.3355"
$ws.Range("E51").Value = "  +1.10%  "
